# This edit corrects a number of rows where the Home/Away match rows had
# their data (id, teams, scores, odds, etc.) swapped / mis-ordered within
# a group of matches that share the same Date. Column A (the running
# index) is left untouched; columns B through AC are rotated among the
# rows listed in each group below so that each row ends up with the data
# that originally belonged to the "next" row in the group (wrapping
# around at the end of the group).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a cyclic group of worksheet row numbers. For a group
# [r0, r1, ..., rN-1], row r(i) receives the B:AC content that originally
# belonged to row r(i+1) (with the last row wrapping back to the first).
$groups = @(
    @(26, 27),
    @(45, 46),
    @(66, 67),
    @(94, 95),
    @(97, 98),
    @(212, 213),
    @(220, 221),
    @(228, 229),
    @(230, 232, 231),
    @(233, 236, 235, 237, 234)
)

foreach ($group in $groups) {
    $count = $group.Length

    # Snapshot the current B:AC values for every row in this group before
    # writing anything back (so reads never see already-modified data).
    $snapshots = @()
    foreach ($rowNum in $group) {
        $rng = $ws.Range("B$rowNum`:AC$rowNum")
        $snapshots += , ($rng.Value2)
    }

    for ($i = 0; $i -lt $count; $i++) {
        $rowNum = $group[$i]
        $sourceIndex = ($i + 1) % $count
        $rng = $ws.Range("B$rowNum`:AC$rowNum")
        $rng.Value2 = $snapshots[$sourceIndex]
    }
}
